# Update "想去人数" (interested-count) figures on both the "展览" and
# "全部类型" worksheets, which contain duplicated data for the same events.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 3427
    $ws.Range("F3").Value = 23
    $ws.Range("F5").Value = 1669
    $ws.Range("F6").Value = 81
}
